# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Thu Sep 19 09:56:31 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.303.17'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').Value = '2.430.39'
$ws.Range('E3').Value = '  +5.31%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.98%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('D9').Value = '2.428.35'
$ws.Range('E9').Value = '  +5.35%  '
$ws.Range('E10').Value = '  +3.77%  '
$ws.Range('E11').Value = '  +4.59%  '
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('E13').Value = '  +5.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +11.98%  '
$ws.Range('D15').Value = '2.864.42'
$ws.Range('D16').Value = '62.193.31'
$ws.Range('E16').Value = '  +3.95%  '
$ws.Range('D18').Value = '2.440.27'
$ws.Range('E18').Value = '  +5.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '347.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.43%  '
$ws.Range('E21').Value = '  +3.11%  '
$ws.Range('E22').Value = '  +3.89%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  +2.65%  '
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  +14.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.10%  '
$ws.Range('E29').Value = '  +14.66%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0786'
$ws.Range('E30').Value = '  +9.17%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.93%  '
$ws.Range('E32').Value = '  +9.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '171.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.36%  '
$ws.Range('E34').Value = '  +5.56%  '
$ws.Range('E35').Value = '  +4.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.58'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.38%  '
$ws.Range('E37').Value = '  +12.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '366.09'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +16.58%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E41').Value = '  +11.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '145.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.33%  '
$ws.Range('E44').Value = '  +7.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.08%  '
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('E47').Value = '  +5.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0516'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.80%  '
$ws.Range('E49').Value = '  +5.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.93%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.70'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +13.90%  '
